$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value = 900
$ws.Cells.Item(111, 9).Value = 900
$ws.Cells.Item(111, 11).Value = 2700
$ws.Cells.Item(111, 13).Value = 367
$ws.Cells.Item(132, 8).Value = 1203.1428
$ws.Cells.Item(132, 9).Value = 1170.3334
$ws.Cells.Item(132, 11).Value = 3511.0002
$ws.Cells.Item(132, 13).Value = -981.0001999999999
$ws.Cells.Item(135, 8).Value = 3249.6667
$ws.Cells.Item(135, 9).Value = 874.5
$ws.Cells.Item(135, 10).Value = 8000
$ws.Cells.Item(135, 11).Value = 7870.5
$ws.Cells.Item(135, 12).Value = 72000
$ws.Cells.Item(135, 13).Value = -5335.5
$ws.Cells.Item(135, 14).Value = -77070
$ws.Cells.Item(137, 8).Value = 3611.2222
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 14).Value = ""
$ws.Cells.Item(141, 8).Value = 2770.5527
$ws.Cells.Item(141, 9).Value = 2025.5454
$ws.Cells.Item(141, 11).Value = 6076.6362
$ws.Cells.Item(141, 13).Value = -896.6361999999999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2732.75
$ws.Cells.Item(61, 9).Value = 2554.9
$ws.Cells.Item(61, 10).Value = 3622
$ws.Cells.Item(61, 11).Value = 2554.9
$ws.Cells.Item(61, 12).Value = 3622
$ws.Cells.Item(61, 13).Value = -2342.9
$ws.Cells.Item(61, 14).Value = -4046
$ws.Cells.Item(74, 8).Value = 3053.125
$ws.Cells.Item(74, 10).Value = 4995
$ws.Cells.Item(74, 12).Value = 4995
$ws.Cells.Item(74, 14).Value = -6743
$ws.Cells.Item(77, 8).Value = 3053.125
$ws.Cells.Item(77, 10).Value = 4995
$ws.Cells.Item(77, 12).Value = 24975
$ws.Cells.Item(77, 14).Value = -33711
$ws.Cells.Item(107, 8).Value = 75000
$ws.Cells.Item(107, 10).Value = 75000
$ws.Cells.Item(107, 12).Value = 75000
$ws.Cells.Item(107, 14).Value = -82680
$ws.Cells.Item(122, 8).Value = 2744.4546
$ws.Cells.Item(122, 9).Value = 2744.4546
$ws.Cells.Item(122, 11).Value = 8233.363799999999
$ws.Cells.Item(122, 13).Value = -5783.363799999999
$ws.Cells.Item(132, 8).Value = 933.6667
$ws.Cells.Item(132, 9).Value = 933.6667
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 2801.0001
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = ""
$ws.Cells.Item(132, 14).Value = -271.0001000000002
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 12).Value = ""
$ws.Cells.Item(133, 14).Value = 0
$ws.Cells.Item(136, 8).Value = 2732.75
$ws.Cells.Item(136, 9).Value = 2554.9
$ws.Cells.Item(136, 10).Value = 3622
$ws.Cells.Item(136, 11).Value = 7664.700000000001
$ws.Cells.Item(136, 12).Value = 10866
$ws.Cells.Item(136, 13).Value = -5114.700000000001
$ws.Cells.Item(136, 14).Value = -15966

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(19, 8).Value = 23903.455
$ws.Cells.Item(19, 9).Value = 24242.5
$ws.Cells.Item(19, 11).Value = 24242.5
$ws.Cells.Item(19, 13).Value = -24069.5
$ws.Cells.Item(107, 8).Value = 411
$ws.Cells.Item(107, 10).Value = 411
$ws.Cells.Item(107, 12).Value = 411
$ws.Cells.Item(107, 14).Value = -4251
$ws.Cells.Item(134, 8).Value = 6214.227
$ws.Cells.Item(134, 9).Value = 6710.2354
$ws.Cells.Item(134, 10).Value = 4527.8
$ws.Cells.Item(134, 11).Value = 20130.7062
$ws.Cells.Item(134, 12).Value = 13583.4
$ws.Cells.Item(134, 13).Value = -17595.7062
$ws.Cells.Item(134, 14).Value = -18653.4

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 1890.3334
$ws.Cells.Item(2, 9).Value = 885
$ws.Cells.Item(2, 10).Value = 2393
$ws.Cells.Item(2, 11).Value = 885
$ws.Cells.Item(2, 12).Value = 2393
$ws.Cells.Item(2, 13).Value = -772
$ws.Cells.Item(2, 14).Value = -2619
$ws.Cells.Item(16, 8).Value = 1670
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 14).Value = ""
$ws.Cells.Item(113, 8).Value = 1670
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 14).Value = ""
$ws.Cells.Item(132, 8).Value = 5196.75
$ws.Cells.Item(132, 9).Value = 5366.3335
$ws.Cells.Item(132, 10).Value = 5095
$ws.Cells.Item(132, 11).Value = 16099.0005
$ws.Cells.Item(132, 12).Value = 15285
$ws.Cells.Item(132, 13).Value = -13569.0005
$ws.Cells.Item(132, 14).Value = -20345

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 685
$ws.Cells.Item(2, 9).Value = 650
$ws.Cells.Item(2, 11).Value = 3900
$ws.Cells.Item(2, 13).Value = -3787
$ws.Cells.Item(34, 8).Value = 1449.5385
$ws.Cells.Item(34, 10).Value = 1893.6666
$ws.Cells.Item(34, 12).Value = 5680.9998
$ws.Cells.Item(34, 14).Value = -5848.9998
$ws.Cells.Item(39, 8).Value = 2277.7778
$ws.Cells.Item(39, 10).Value = 2277.7778
$ws.Cells.Item(39, 12).Value = 6833.3334
$ws.Cells.Item(39, 14).Value = -7421.3334
$ws.Cells.Item(55, 8).Value = 420.5
$ws.Cells.Item(55, 10).Value = 735
$ws.Cells.Item(55, 12).Value = 2205
$ws.Cells.Item(55, 14).Value = -2559
$ws.Cells.Item(68, 8).Value = 1299.8
$ws.Cells.Item(68, 9).Value = 999.5
$ws.Cells.Item(68, 11).Value = 2998.5
$ws.Cells.Item(68, 13).Value = -2187.5
$ws.Cells.Item(71, 8).Value = 1299.8
$ws.Cells.Item(71, 9).Value = 999.5
$ws.Cells.Item(71, 11).Value = 8995.5
$ws.Cells.Item(71, 13).Value = -4939.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 300
$ws.Cells.Item(107, 9).Value = 300
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 300
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = ""
$ws.Cells.Item(107, 14).Value = 1620
$ws.Cells.Item(113, 8).Value = 1155
$ws.Cells.Item(113, 9).Value = 1014.1667
$ws.Cells.Item(113, 11).Value = 1014.1667
$ws.Cells.Item(113, 13).Value = 1155.8333
$ws.Cells.Item(132, 8).Value = 2499
$ws.Cells.Item(132, 9).Value = 2998
$ws.Cells.Item(132, 10).Value = 2000
$ws.Cells.Item(132, 11).Value = 8994
$ws.Cells.Item(132, 12).Value = 6000
$ws.Cells.Item(132, 13).Value = -6464
$ws.Cells.Item(132, 14).Value = -11060

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 22503.908
$ws.Cells.Item(7, 9).Value = 22503.908
$ws.Cells.Item(7, 11).Value = 22503.908
$ws.Cells.Item(7, 13).Value = -22391.908
$ws.Cells.Item(17, 8).Value = 21329
$ws.Cells.Item(17, 9).Value = 12000
$ws.Cells.Item(17, 10).Value = 25993.5
$ws.Cells.Item(17, 11).Value = 12000
$ws.Cells.Item(17, 12).Value = 25993.5
$ws.Cells.Item(17, 13).Value = -11830
$ws.Cells.Item(17, 14).Value = -26333.5
$ws.Cells.Item(122, 8).Value = 3152
$ws.Cells.Item(122, 9).Value = 3152
$ws.Cells.Item(122, 11).Value = 9456
$ws.Cells.Item(122, 13).Value = -7006
$ws.Cells.Item(126, 8).Value = 22503.908
$ws.Cells.Item(126, 9).Value = 22503.908
$ws.Cells.Item(126, 11).Value = 67511.724
$ws.Cells.Item(126, 13).Value = -65041.724
$ws.Cells.Item(132, 8).Value = 2975
$ws.Cells.Item(132, 9).Value = 2975
$ws.Cells.Item(132, 11).Value = 8925
$ws.Cells.Item(132, 13).Value = -6395
$ws.Cells.Item(136, 8).Value = 3972.9285
$ws.Cells.Item(136, 9).Value = 3067.1428
$ws.Cells.Item(136, 11).Value = 9201.428400000001
$ws.Cells.Item(136, 13).Value = -6651.428400000001

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(42, 8).Value = 10000
$ws.Cells.Item(42, 9).Value = 10000
$ws.Cells.Item(42, 11).Value = 10000
$ws.Cells.Item(42, 13).Value = -9622
$ws.Cells.Item(43, 8).Value = 20000
$ws.Cells.Item(43, 9).Value = 20000
$ws.Cells.Item(43, 11).Value = 20000
$ws.Cells.Item(43, 13).Value = -19851
$ws.Cells.Item(122, 8).Value = 2233.9
$ws.Cells.Item(122, 9).Value = 2234.1428
$ws.Cells.Item(122, 10).Value = 2233.3333
$ws.Cells.Item(122, 11).Value = 6702.428400000001
$ws.Cells.Item(122, 12).Value = 6699.999899999999
$ws.Cells.Item(122, 13).Value = -4252.428400000001
$ws.Cells.Item(122, 14).Value = -11599.9999
$ws.Cells.Item(136, 8).Value = 8312.923000000001
$ws.Cells.Item(136, 10).Value = 9224.286
$ws.Cells.Item(136, 12).Value = 27672.858
$ws.Cells.Item(136, 13).Value = -32772.858
